$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H98").Value = 6187.5
$ws.Range("I98").Value = 3742.8572
$ws.Range("J98").Value = 8088.8887
$ws.Range("K98").Value = 3742.8572
$ws.Range("L98").Value = 8088.8887
$ws.Range("M98").Value = -2244.8572
$ws.Range("N98").Value = -11084.8887

$ws.Range("H100").Value = 40002468
$ws.Range("I100").Value = 40002468
$ws.Range("K100").Value = 40002468
$ws.Range("M100").Value = -40001927

$ws.Range("H106").Value = 4334.091
$ws.Range("I106").Value = 1812.5
$ws.Range("J106").Value = 7360
$ws.Range("K106").Value = 1812.5
$ws.Range("L106").Value = 7360
$ws.Range("M106").Value = -1181.5
$ws.Range("N106").Value = -8622

$ws.Range("H113").Value = 5092
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5092
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5092
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11600

$ws.Range("H122").Value = 6187.5
$ws.Range("I122").Value = 3742.8572
$ws.Range("J122").Value = 8088.8887
$ws.Range("K122").Value = 11228.5716
$ws.Range("L122").Value = 24266.6661
$ws.Range("M122").Value = -8778.571599999999
$ws.Range("N122").Value = -29166.6661

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 3863.7031
$ws.Range("I32").Value = 4551.114
$ws.Range("J32").Value = 2351.4
$ws.Range("K32").Value = 4551.114
$ws.Range("L32").Value = 2351.4
$ws.Range("M32").Value = -4264.114
$ws.Range("N32").Value = -2925.4

$ws.Range("H45").Value = 2013.3334
$ws.Range("I45").Value = 1816
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1816
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1439
$ws.Range("N45").Value = -3754

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 3455
$ws.Range("I31").Value = 1211.0834
$ws.Range("J31").Value = 5526.3076
$ws.Range("K31").Value = 1211.0834
$ws.Range("L31").Value = 5526.3076
$ws.Range("M31").Value = -916.0834
$ws.Range("N31").Value = -6116.3076

$ws.Range("H34").Value = 3455
$ws.Range("I34").Value = 1211.0834
$ws.Range("J34").Value = 5526.3076
$ws.Range("K34").Value = 1211.0834
$ws.Range("L34").Value = 5526.3076
$ws.Range("M34").Value = -1009.0834
$ws.Range("N34").Value = -5930.3076

$ws.Range("H99").Value = 14290714
$ws.Range("J99").Value = 6588.8887
$ws.Range("L99").Value = 6588.8887
$ws.Range("N99").Value = -9584.8887

$ws.Range("H122").Value = 2338.6155
$ws.Range("I122").Value = 1140.3
$ws.Range("J122").Value = 6333
$ws.Range("K122").Value = 3420.9
$ws.Range("L122").Value = 18999
$ws.Range("M122").Value = -970.8999999999996
$ws.Range("N122").Value = -23899

$ws.Range("H126").Value = 14290714
$ws.Range("J126").Value = 6588.8887
$ws.Range("L126").Value = 19766.6661
$ws.Range("N126").Value = -24706.6661

$ws.Range("H132").Value = 3537.2222
$ws.Range("I132").Value = 3079.55
$ws.Range("J132").Value = 4844.857
$ws.Range("K132").Value = 9238.650000000001
$ws.Range("L132").Value = 14534.571
$ws.Range("M132").Value = -6708.650000000001
$ws.Range("N132").Value = -19594.571

$ws.Range("H134").Value = 8888.134
$ws.Range("I134").Value = 12266.223
$ws.Range("J134").Value = 3821
$ws.Range("K134").Value = 36798.669
$ws.Range("L134").Value = 11463
$ws.Range("M134").Value = -34263.669
$ws.Range("N134").Value = -16533

$ws.Range("H137").Value = 45935.715
$ws.Range("J137").Value = 45935.715
$ws.Range("L137").Value = 45935.715
$ws.Range("N137").Value = -56135.715

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 668423.75
$ws.Range("I5").Value = 556.1
$ws.Range("J5").Value = 1336291.4
$ws.Range("K5").Value = 1668.3
$ws.Range("L5").Value = 4008874.2
$ws.Range("M5").Value = -1556.3
$ws.Range("N5").Value = -4009098.2

$ws.Range("H56").Value = 4411.5386
$ws.Range("I56").Value = 4411.5386
$ws.Range("K56").Value = 4411.5386
$ws.Range("M56").Value = -3881.5386

$ws.Range("H129").Value = 3234.3845
$ws.Range("I129").Value = 4250
$ws.Range("J129").Value = 2599.625
$ws.Range("K129").Value = 12750
$ws.Range("L129").Value = 7798.875
$ws.Range("M129").Value = -7750
$ws.Range("N129").Value = -17798.875

$ws.Range("H135").Value = 668423.75
$ws.Range("I135").Value = 556.1
$ws.Range("J135").Value = 1336291.4
$ws.Range("K135").Value = 5004.900000000001
$ws.Range("L135").Value = 12026622.6
$ws.Range("M135").Value = -2469.900000000001
$ws.Range("N135").Value = -12031692.6

$ws.Range("H137").Value = 2503.25
$ws.Range("I137").Value = 1490
$ws.Range("J137").Value = 3516.5
$ws.Range("K137").Value = 4470
$ws.Range("L137").Value = 10549.5
$ws.Range("M137").Value = 630
$ws.Range("N137").Value = -20749.5

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 6453.731
$ws.Range("I70").Value = 5819.857
$ws.Range("K70").Value = 5819.857
$ws.Range("M70").Value = -5549.857

$ws.Range("H73").Value = 6453.731
$ws.Range("I73").Value = 5819.857
$ws.Range("K73").Value = 5819.857
$ws.Range("M73").Value = -4883.857

$ws.Range("H102").Value = 3028.8125
$ws.Range("I102").Value = 2038.5
$ws.Range("J102").Value = 5999.75
$ws.Range("K102").Value = 2038.5
$ws.Range("L102").Value = 5999.75
$ws.Range("M102").Value = -416.5
$ws.Range("N102").Value = -9243.75

$ws.Range("H122").Value = 4598.4614
$ws.Range("I122").Value = 2252.7273
$ws.Range("K122").Value = 6758.1819
$ws.Range("M122").Value = -4308.1819

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 4747.067
$ws.Range("I7").Value = 3766.375
$ws.Range("J7").Value = 5867.857
$ws.Range("K7").Value = 3766.375
$ws.Range("L7").Value = 5867.857
$ws.Range("M7").Value = -3654.375
$ws.Range("N7").Value = -6091.857

$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9864

$ws.Range("H122").Value = 5045
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 5257.0835
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 15771.2505
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -20671.2505

$ws.Range("H126").Value = 4747.067
$ws.Range("I126").Value = 3766.375
$ws.Range("J126").Value = 5867.857
$ws.Range("K126").Value = 11299.125
$ws.Range("L126").Value = 17603.571
$ws.Range("M126").Value = -8829.125
$ws.Range("N126").Value = -22543.571

$ws.Range("H132").Value = 22736.375
$ws.Range("I132").Value = 240000
$ws.Range("J132").Value = 8252.134
$ws.Range("K132").Value = 720000
$ws.Range("L132").Value = 24756.402
$ws.Range("M132").Value = -717470
$ws.Range("N132").Value = -29816.402

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H108").Value = 36375.332
$ws.Range("J108").Value = 36375.332
$ws.Range("L108").Value = 36375.332
$ws.Range("N108").Value = -44055.332

$ws.Range("H132").Value = 14496349
$ws.Range("I132").Value = 2197.7144
$ws.Range("J132").Value = 20837540
$ws.Range("K132").Value = 6593.1432
$ws.Range("L132").Value = 62512620
$ws.Range("M132").Value = -4063.1432
$ws.Range("N132").Value = -62517680
